$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.775841999999999
$ws.Range("H2").Value = 8.327525999999999
$ws.Range("I2").Value = 0.0624750527258915
$ws.Range("J2").Value = 0.0624750527258915
$ws.Range("M2").Value = 0.6415476666666667
$ws.Range("N2").Value = 1.924643
$ws.Range("O2").Value = 0.1426849042655057
$ws.Range("P2").Value = 0.1426849042655057
$ws.Range("Q2").Value = 1.780834958135333
$ws.Range("R2").Value = 16.027514623218
$ws.Range("S2").Value = 0.00891424691717625
$ws.Range("T2").Value = 0.008914246917176252
$ws.Range("G3").Value = 2.775841999999999
$ws.Range("H3").Value = 8.327525999999999
$ws.Range("I3").Value = 0.0624750527258915
$ws.Range("J3").Value = 0.0624750527258915
$ws.Range("O3").Value = 0.5986102210699216
$ws.Range("P3").Value = 0.5986102210699217
$ws.Range("Q3").Value = 7.471189846367999
$ws.Range("R3").Value = 67.240708617312
$ws.Range("S3").Value = 0.03739820512360091
$ws.Range("T3").Value = 0.03739820512360092
$ws.Range("G4").Value = 2.775841999999999
$ws.Range("H4").Value = 8.327525999999999
$ws.Range("I4").Value = 0.0624750527258915
$ws.Range("J4").Value = 0.0624750527258915
$ws.Range("M4").Value = 1.163203
$ws.Range("N4").Value = 3.489609
$ws.Range("O4").Value = 0.2587048746645726
$ws.Range("P4").Value = 0.2587048746645726
$ws.Range("Q4").Value = 3.228867741925999
$ws.Range("R4").Value = 29.05980967733399
$ws.Range("S4").Value = 0.01616260068511433
$ws.Range("T4").Value = 0.01616260068511433
$ws.Range("I5").Value = 0.2652892219050753
$ws.Range("J5").Value = 0.2652892219050753
$ws.Range("M5").Value = 0.6415476666666667
$ws.Range("N5").Value = 1.924643
$ws.Range("O5").Value = 0.1426849042655057
$ws.Range("P5").Value = 0.1426849042655057
$ws.Range("Q5").Value = 7.561999546569224
$ws.Range("R5").Value = 68.057995919123
$ws.Range("S5").Value = 0.03785276723019617
$ws.Range("T5").Value = 0.03785276723019617
$ws.Range("I6").Value = 0.2652892219050753
$ws.Range("J6").Value = 0.2652892219050753
$ws.Range("O6").Value = 0.5986102210699216
$ws.Range("P6").Value = 0.5986102210699217
$ws.Range("S6").Value = 0.1588048397720646
$ws.Range("T6").Value = 0.1588048397720646
$ws.Range("I7").Value = 0.2652892219050753
$ws.Range("J7").Value = 0.2652892219050753
$ws.Range("M7").Value = 1.163203
$ws.Range("N7").Value = 3.489609
$ws.Range("O7").Value = 0.2587048746645726
$ws.Range("P7").Value = 0.2587048746645726
$ws.Range("Q7").Value = 13.71081373309433
$ws.Range("R7").Value = 123.397323597849
$ws.Range("S7").Value = 0.0686316149028145
$ws.Range("T7").Value = 0.0686316149028145
$ws.Range("G8").Value = 29.86824466666667
$ws.Range("H8").Value = 89.60473400000001
$ws.Range("I8").Value = 0.6722357253690333
$ws.Range("J8").Value = 0.6722357253690333
$ws.Range("M8").Value = 0.6415476666666667
$ws.Range("N8").Value = 1.924643
$ws.Range("O8").Value = 0.1426849042655057
$ws.Range("P8").Value = 0.1426849042655057
$ws.Range("Q8").Value = 19.16190267332912
$ws.Range("R8").Value = 172.457124059962
$ws.Range("S8").Value = 0.09591789011813331
$ws.Range("T8").Value = 0.09591789011813331
$ws.Range("G9").Value = 29.86824466666667
$ws.Range("H9").Value = 89.60473400000001
$ws.Range("I9").Value = 0.6722357253690333
$ws.Range("J9").Value = 0.6722357253690333
$ws.Range("O9").Value = 0.5986102210699216
$ws.Range("P9").Value = 0.5986102210699217
$ws.Range("Q9").Value = 80.390499993312
$ws.Range("R9").Value = 723.514499939808
$ws.Range("S9").Value = 0.4024071761742561
$ws.Range("T9").Value = 0.4024071761742562
$ws.Range("G10").Value = 29.86824466666667
$ws.Range("H10").Value = 89.60473400000001
$ws.Range("I10").Value = 0.6722357253690333
$ws.Range("J10").Value = 0.6722357253690333
$ws.Range("M10").Value = 1.163203
$ws.Range("N10").Value = 3.489609
$ws.Range("O10").Value = 0.2587048746645726
$ws.Range("P10").Value = 0.2587048746645726
$ws.Range("Q10").Value = 34.74283180100067
$ws.Range("R10").Value = 312.685486209006
$ws.Range("S10").Value = 0.1739106590766439
$ws.Range("T10").Value = 0.1739106590766439
